# Generate Report for Handback
#
# This script updates the localization-status workbook after a handback run:
#  - Marks the two content items in the "zh-cn" and "de-de" sheets as handed
#    back (Status column).
#  - Fills in the "Latest Target File" (hyperlink to the source .md on
#    GitHub) and "Latest Handback File" (generated .xlf file name) columns.
#  - Stamps the "Latest Handback DateTime" for each language sheet.
#  - Widens a few columns so the new/longer text is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# GitHub source-file URLs (identical in every language sheet); reused as the
# hyperlink target for the new "Latest Target File" column.
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/375a2a5653584246c87b6f395ae6645bfb86fa99/e2e/0bfa8ba1-5a98-4a68-a520-6aa45212fba1.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/375a2a5653584246c87b6f395ae6645bfb86fa99/e2e/d67d93d3-4252-4765-8e33-fa077ba03446.md"
$mdName1 = "0bfa8ba1-5a98-4a68-a520-6aa45212fba1.md"
$mdName2 = "d67d93d3-4252-4765-8e33-fa077ba03446.md"

# Per-language sheet configuration: handback file names and handback time.
$langInfo = @{
    "zh-cn" = @{
        HandbackFile1 = "0bfa8ba1-5a98-4a68-a520-6aa45212fba1.e0d28628a5970ca168296536be8c126696277144.zh-cn.xlf"
        HandbackFile2 = "d67d93d3-4252-4765-8e33-fa077ba03446.04bfbcc5ecd528731ae4bba5eb806ddbb71dfc58.zh-cn.xlf"
        HandbackTime  = "2016-08-19 06:29:28"
    }
    "de-de" = @{
        HandbackFile1 = "0bfa8ba1-5a98-4a68-a520-6aa45212fba1.e0d28628a5970ca168296536be8c126696277144.de-de.xlf"
        HandbackFile2 = "d67d93d3-4252-4765-8e33-fa077ba03446.04bfbcc5ecd528731ae4bba5eb806ddbb71dfc58.de-de.xlf"
        HandbackTime  = "2016-08-19 06:29:35"
    }
}

foreach ($langName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($langName)
    $info = $langInfo[$langName]

    # Status column (C) -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (I) -> hyperlink to the source markdown file.
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdName2)

    # Latest Handback File (J) -> generated xliff file name.
    $ws.Range("J2").Value = $info.HandbackFile1
    $ws.Range("J3").Value = $info.HandbackFile2

    # Latest Handback DateTime (K) -> timestamp of this handback run.
    $ws.Range("K2").Value = $info.HandbackTime
    $ws.Range("K3").Value = $info.HandbackTime

    # Widen columns so the new content fits comfortably.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

# Overview sheet: widen the zh-cn / de-de status columns (E, F) to match.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668
